# Auto-generated Excel COM-interop script
# Applies the "Updated cryptos list" diff: updates Price (D) / Volume(1h) (E)
# figures for rows 2-51, and refreshes the Coin name/Link/Price/Volume for
# rows 40-42 where the ranking order shifted (Fetch.AI, dogwifhat, Stacks).
#
# D/E columns hold numeric-looking values stored as TEXT (e.g. "70.520.83",
# "  +0.01%  "). Plain `.Value = "2.94"` would get auto-coerced to a real
# number by Excel, so we prefix with a literal apostrophe to force text,
# then reset `.Style` back to "Normal" so the transient quote-prefix style
# Excel applies does not linger as a visible formatting change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.520.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.01%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.619.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'584.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.54%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.616.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.18%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.633"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.75%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +4.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.666"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'56.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +7.50%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.62%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.200.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.02%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'19.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.37%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.615.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'70.518.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.09%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.19%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.56%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'494.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.16%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'19.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.82%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'4.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -7.15%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'97.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +7.04%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.95%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.27%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'32.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.98%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.95%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'66.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'584.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -8.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'39.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.83%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0₃0823"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.17%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.37%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.64%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +19.73%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.59%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.244.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -6.59%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.98%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0448"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +7.19%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.64%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.139"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'3.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.26%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.04%  "
$ws.Range("E51").Style = "Normal"
